# Updated legacy GSC export data.
#
# The workbook tracks a rolling window of daily video-indexing stats on the
# "Chart" sheet (A:Date, B:No video indexed, C:Video indexed, D:Impressions).
# The export was refreshed: the two oldest dates (2025-11-21, 2025-11-22)
# dropped out of the window and every remaining row shifted up by two, with
# no new rows appended at the bottom (the sheet shrinks from 89 to 87 rows).
#
# Deleting the two oldest rows (rather than re-typing date strings into
# cells) lets Excel move the existing text/number cells as-is, so the
# date strings stay text (no auto-conversion to date serials) and cell
# styles/number formats are untouched, exactly like the real edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
